$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.7471093333333333
$ws.Range("H2").Value = 2.241328
$ws.Range("I2").Value = 0.4648531405557055
$ws.Range("J2").Value = 0.4648531405557054
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.446352
$ws.Range("N2").Value = 4.339056
$ws.Range("O2").Value = 0.05318694539780245
$ws.Range("P2").Value = 0.05318694539780246
$ws.Range("Q2").Value = 1.080583078485333
$ws.Range("R2").Value = 9.725247706368
$ws.Range("S2").Value = 0.02472411860473329
$ws.Range("T2").Value = 0.02472411860473329

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.7471093333333333
$ws.Range("H3").Value = 2.241328
$ws.Range("I3").Value = 0.4648531405557055
$ws.Range("J3").Value = 0.4648531405557054
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.619945333333334
$ws.Range("N3").Value = 16.859836
$ws.Range("O3").Value = 0.2066631951161506
$ws.Range("P3").Value = 0.2066631951161506
$ws.Range("Q3").Value = 4.198713611356444
$ws.Range("R3").Value = 37.788422502208
$ws.Range("S3").Value = 0.09606803528701915
$ws.Range("T3").Value = 0.09606803528701914

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.7471093333333333
$ws.Range("H4").Value = 2.241328
$ws.Range("I4").Value = 0.4648531405557055
$ws.Range("J4").Value = 0.4648531405557054
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.371955333333334
$ws.Range("N4").Value = 7.115866
$ws.Range("O4").Value = 0.08722431247720218
$ws.Range("P4").Value = 0.08722431247720218
$ws.Range("Q4").Value = 1.772109967783111
$ws.Range("R4").Value = 15.948989710048
$ws.Range("S4").Value = 0.04054649558783964
$ws.Range("T4").Value = 0.04054649558783963

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.7471093333333333
$ws.Range("H5").Value = 2.241328
$ws.Range("I5").Value = 0.4648531405557055
$ws.Range("J5").Value = 0.4648531405557054
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.306061
$ws.Range("N5").Value = 15.918183
$ws.Range("O5").Value = 0.1951206737256277
$ws.Range("P5").Value = 0.1951206737256277
$ws.Range("Q5").Value = 3.964207696336
$ws.Range("R5").Value = 35.67786926702399
$ws.Range("S5").Value = 0.09070245796870316
$ws.Range("T5").Value = 0.09070245796870313

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.7471093333333333
$ws.Range("H6").Value = 2.241328
$ws.Range("I6").Value = 0.4648531405557055
$ws.Range("J6").Value = 0.4648531405557054
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.866195333333334
$ws.Range("N6").Value = 29.598586
$ws.Range("O6").Value = 0.3628112606599592
$ws.Range("P6").Value = 0.3628112606599592
$ws.Range("Q6").Value = 7.371126618023111
$ws.Range("R6").Value = 66.34013956220799
$ws.Range("S6").Value = 0.1686539539467567
$ws.Range("T6").Value = 0.1686539539467567

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.7471093333333333
$ws.Range("H7").Value = 2.241328
$ws.Range("I7").Value = 0.4648531405557055
$ws.Range("J7").Value = 0.4648531405557054
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.583231666666667
$ws.Range("N7").Value = 7.749695
$ws.Range("O7").Value = 0.09499361262325784
$ws.Range("P7").Value = 0.09499361262325784
$ws.Range("Q7").Value = 1.929956488328889
$ws.Range("R7").Value = 17.36960839496
$ws.Range("S7").Value = 0.04415807916065351
$ws.Range("T7").Value = 0.04415807916065351

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4799286666666667
$ws.Range("H8").Value = 1.439786
$ws.Range("I8").Value = 0.2986127170267525
$ws.Range("J8").Value = 0.2986127170267524
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.446352
$ws.Range("N8").Value = 4.339056
$ws.Range("O8").Value = 0.05318694539780245
$ws.Range("P8").Value = 0.05318694539780246
$ws.Range("Q8").Value = 0.6941457868906667
$ws.Range("R8").Value = 6.247312082016
$ws.Range("S8").Value = 0.01588229827559132
$ws.Range("T8").Value = 0.01588229827559132

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4799286666666667
$ws.Range("H9").Value = 1.439786
$ws.Range("I9").Value = 0.2986127170267525
$ws.Range("J9").Value = 0.2986127170267524
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.619945333333334
$ws.Range("N9").Value = 16.859836
$ws.Range("O9").Value = 0.2066631951161506
$ws.Range("P9").Value = 0.2066631951161506
$ws.Range("Q9").Value = 2.697172870566222
$ws.Range("R9").Value = 24.274555835096
$ws.Range("S9").Value = 0.06171225820306362
$ws.Range("T9").Value = 0.06171225820306361

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4799286666666667
$ws.Range("H10").Value = 1.439786
$ws.Range("I10").Value = 0.2986127170267525
$ws.Range("J10").Value = 0.2986127170267524
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.371955333333334
$ws.Range("N10").Value = 7.115866
$ws.Range("O10").Value = 0.08722431247720218
$ws.Range("P10").Value = 0.08722431247720218
$ws.Range("Q10").Value = 1.138369360519556
$ws.Range("R10").Value = 10.245324244676
$ws.Range("S10").Value = 0.02604628893960781
$ws.Range("T10").Value = 0.0260462889396078

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.4799286666666667
$ws.Range("H11").Value = 1.439786
$ws.Range("I11").Value = 0.2986127170267525
$ws.Range("J11").Value = 0.2986127170267524
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.306061
$ws.Range("N11").Value = 15.918183
$ws.Range("O11").Value = 0.1951206737256277
$ws.Range("P11").Value = 0.1951206737256277
$ws.Range("Q11").Value = 2.546530780982
$ws.Range("R11").Value = 22.918777028838
$ws.Range("S11").Value = 0.05826551452930016
$ws.Range("T11").Value = 0.05826551452930014

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.4799286666666667
$ws.Range("H12").Value = 1.439786
$ws.Range("I12").Value = 0.2986127170267525
$ws.Range("J12").Value = 0.2986127170267524
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 9.866195333333334
$ws.Range("N12").Value = 29.598586
$ws.Range("O12").Value = 0.3628112606599592
$ws.Range("P12").Value = 0.3628112606599592
$ws.Range("Q12").Value = 4.735069971399556
$ws.Range("R12").Value = 42.615629742596
$ws.Range("S12").Value = 0.1083400563135717
$ws.Range("T12").Value = 0.1083400563135717

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.4799286666666667
$ws.Range("H13").Value = 1.439786
$ws.Range("I13").Value = 0.2986127170267525
$ws.Range("J13").Value = 0.2986127170267524
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.583231666666667
$ws.Range("N13").Value = 7.749695
$ws.Range("O13").Value = 0.09499361262325784
$ws.Range("P13").Value = 0.09499361262325784
$ws.Range("Q13").Value = 1.239766929474444
$ws.Range("R13").Value = 11.15790236527
$ws.Range("S13").Value = 0.02836630076561783
$ws.Range("T13").Value = 0.02836630076561783

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.3801563333333333
$ws.Range("H14").Value = 1.140469
$ws.Range("I14").Value = 0.2365341424175421
$ws.Range("J14").Value = 0.2365341424175421
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.446352
$ws.Range("N14").Value = 4.339056
$ws.Range("O14").Value = 0.05318694539780245
$ws.Range("P14").Value = 0.05318694539780246
$ws.Range("Q14").Value = 0.5498398730293333
$ws.Range("R14").Value = 4.948558857264
$ws.Range("S14").Value = 0.01258052851747784
$ws.Range("T14").Value = 0.01258052851747784

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.3801563333333333
$ws.Range("H15").Value = 1.140469
$ws.Range("I15").Value = 0.2365341424175421
$ws.Range("J15").Value = 0.2365341424175421
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.619945333333334
$ws.Range("N15").Value = 16.859836
$ws.Range("O15").Value = 0.2066631951161506
$ws.Range("P15").Value = 0.2066631951161506
$ws.Range("Q15").Value = 2.136457811453778
$ws.Range("R15").Value = 19.228120303084
$ws.Range("S15").Value = 0.04888290162606788
$ws.Range("T15").Value = 0.04888290162606786

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.3801563333333333
$ws.Range("H16").Value = 1.140469
$ws.Range("I16").Value = 0.2365341424175421
$ws.Range("J16").Value = 0.2365341424175421
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2.371955333333334
$ws.Range("N16").Value = 7.115866
$ws.Range("O16").Value = 0.08722431247720218
$ws.Range("P16").Value = 0.08722431247720218
$ws.Range("Q16").Value = 0.9017138423504445
$ws.Range("R16").Value = 8.115424581154
$ws.Range("S16").Value = 0.02063152794975474
$ws.Range("T16").Value = 0.02063152794975473

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.3801563333333333
$ws.Range("H17").Value = 1.140469
$ws.Range("I17").Value = 0.2365341424175421
$ws.Range("J17").Value = 0.2365341424175421
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.306061
$ws.Range("N17").Value = 15.918183
$ws.Range("O17").Value = 0.1951206737256277
$ws.Range("P17").Value = 0.1951206737256277
$ws.Range("Q17").Value = 2.017132694203
$ws.Range("R17").Value = 18.154194247827
$ws.Range("S17").Value = 0.0461527012276244
$ws.Range("T17").Value = 0.04615270122762438

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 0.3801563333333333
$ws.Range("H18").Value = 1.140469
$ws.Range("I18").Value = 0.2365341424175421
$ws.Range("J18").Value = 0.2365341424175421
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 9.866195333333334
$ws.Range("N18").Value = 29.598586
$ws.Range("O18").Value = 0.3628112606599592
$ws.Range("P18").Value = 0.3628112606599592
$ws.Range("Q18").Value = 3.750696641870444
$ws.Range("R18").Value = 33.756269776834
$ws.Range("S18").Value = 0.0858172503996308
$ws.Range("T18").Value = 0.08581725039963078

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 0.3801563333333333
$ws.Range("H19").Value = 1.140469
$ws.Range("I19").Value = 0.2365341424175421
$ws.Range("J19").Value = 0.2365341424175421
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 2.583231666666667
$ws.Range("N19").Value = 7.749695
$ws.Range("O19").Value = 0.09499361262325784
$ws.Range("P19").Value = 0.09499361262325784
$ws.Range("Q19").Value = 0.9820318785505556
$ws.Range("R19").Value = 8.838286906955
$ws.Range("S19").Value = 0.0224692326969865
$ws.Range("T19").Value = 0.0224692326969865
